$d = $word.ActiveDocument

# 1. "Назначение программы" sentence -> new purpose text
$d.Content.Find.Execute(
    "Назначение программы: поиск среди элементов массива наибольшего числа, кратного 4.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Назначение программы: нахождение значения функции:", 2) | Out-Null

# 2. Split off an empty paragraph before the "Вывод:" paragraph, so the
#    page-break marker ends up attached to the (now second) "Вывод:" run.
$vyvodRng = $d.Content
$vyvodRng.Find.Execute("Вывод:", $false) | Out-Null
$vyvodPara = $vyvodRng.Paragraphs(1)
$vyvodPara.Range.InsertParagraphBefore()

# 3. Replace the "Вывод" conclusion paragraph text with the lab7 content.
$oldConclusion = "В процессе выполнения лабораторной работы был получен опыт работы с циклическими программами и работы с одномерными массивами. Были изучены различные виды адресации (прямая абсолютная и относительная, косвенная относительная). Данный материал поможет мне в будущем при написании различного рода программ, использующих циклы и одномерные массивы"
$newConclusion = "В процессе выполнения лабораторной работы был получен опыт работы с подпрограммами и стеком, разобрался, каким образом реализован стек в БЭВМ, изучил принцип действия команд PUSH, POP, CALL, RET"
$d.Content.Find.Execute($oldConclusion, $true, $false, $false, $false, $false, $true, 1, $false, $newConclusion, 2) | Out-Null

# 4. Tag the instruction mnemonics in the new conclusion sentence as en-US,
#    matching how the rest of the document marks inline English text.
#    Restrict the search to the conclusion paragraph only so we don't touch
#    the pre-existing PUSH/POP/CALL/RET occurrences in the instruction table.
foreach ($word_ in @("PUSH", "POP", "CALL", "RET")) {
    $concRng = $d.Content
    $concRng.Find.Execute("принцип действия команд", $false) | Out-Null
    $concPara = $concRng.Paragraphs(1)
    $boundedRng = $d.Range($concPara.Range.Start, $concPara.Range.End)

    $boundedRng.Find.Execute($word_, $true, $false, $false, $false, $false, $true) | Out-Null
    $boundedRng.LanguageID = "en-US"
}
